# Attendance app entrypoint: append today's attendance check-ins to the
# end of the "2024-05-16" sheet (below the existing table, leaving one
# blank separator row, matching how previous day's sheet was written).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-05-16")

# Existing data occupies rows 1-5 (header + 4 students); a blank row 6
# separates it from the freshly appended batch starting at row 7.
$startRow = 7

$records = @(
    @{ Roll = "1978"; Attendance = "Absent";  Timestamp = "None" },
    @{ Roll = "1990"; Attendance = "Absent";  Timestamp = "None" },
    @{ Roll = "2010"; Attendance = "Present"; Timestamp = "2024-05-16 12:38:05" },
    @{ Roll = "2089"; Attendance = "Present"; Timestamp = "2024-05-16 12:38:08" }
)

$r = $startRow
foreach ($rec in $records) {
    # Leading apostrophe forces the roll number to be stored as text
    # (it would otherwise be auto-detected as a number), matching the
    # string type used for the roll numbers already in the sheet.
    $ws.Range("A$r").Value = "'" + $rec.Roll
    $ws.Range("B$r").Value = $rec.Attendance
    $ws.Range("C$r").Value = $rec.Timestamp
    $r++
}
